# [Kadastro App] Yeni kayit eklendi: 1
# Adds the new "Kayıt" (record) row to both the master "Kayitlar" sheet and
# the matching district sheet ("Silifke", since Birim = "Silifke").

$wb = $excel.ActiveWorkbook

$recordNo   = "1"
$tarih      = "2025-09-05"
$birim      = "Silifke"
$parsel     = "5"
$is         = "KAMULAŞTIRMA"
$personeller = "EMİNE ALANLI KIRCILI (K.Mühendisi), AHMET YILDIRIM (K.Teknisyeni)"

$sheetNames = @("Kayitlar", "Silifke")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Sheets.Item($sheetName)

    # Numeric-looking / date-looking values must stay text (the sheet keeps
    # them "numberStoredAsText"), so force text format before assigning.
    $ws.Range("A2").NumberFormat = "@"
    $ws.Range("A2").Value = $recordNo

    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = $tarih

    $ws.Range("C2").Value = $birim

    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = $parsel

    $ws.Range("E2").Value = $is

    $ws.Range("F2").Value = $personeller
}
